# Add a new restaurant entry ("Lake Inez") as row 20 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = "Lake Inez"
$ws.Cells.Item($row, 2).Value = "Leslieville"
$ws.Cells.Item($row, 3).Value = "Secret patio menu, really anything on the menu is going to be outrageous; as of August 2024, our favorite restaurant in Toronto"
$ws.Cells.Item($row, 4).Value = "Farm to Table but also kinda Asian?"
$ws.Cells.Item($row, 5).Value = 43.673155291418801
$ws.Cells.Item($row, 6).Value = -79.3208615612443

# Move the active selection the way Excel leaves it after entering a new
# row and pressing Enter/Tab past the last filled column.
$ws.Range("D21").Select()
